$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1500
$ws.Range("I6").Value = 1500
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 4500
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -4388
# Row 34
$ws.Range("H34").Value = 3900
$ws.Range("I34").Value = 1216.6666
$ws.Range("K34").Value = 1216.6666
$ws.Range("M34").Value = -1013.6666
# Row 36
$ws.Range("H36").Value = 3900
$ws.Range("I36").Value = 1216.6666
$ws.Range("K36").Value = 1216.6666
$ws.Range("M36").Value = -501.6666
# Row 40
$ws.Range("H40").Value = 1528.6765
$ws.Range("I40").Value = 907.4167
$ws.Range("J40").Value = 1867.5454
$ws.Range("K40").Value = 907.4167
$ws.Range("L40").Value = 1867.5454
$ws.Range("M40").Value = -732.4167
$ws.Range("N40").Value = -2217.5454
# Row 137
$ws.Range("H137").Value = 32189.666
$ws.Range("I137").Value = 1616.8334
$ws.Range("J137").Value = 93335.336
$ws.Range("K137").Value = 4850.5002
$ws.Range("L137").Value = 280006.008
$ws.Range("M137").Value = -2300.5002
$ws.Range("N137").Value = -285106.008
# Row 138
$ws.Range("H138").Value = 5468321
$ws.Range("I138").Value = 10754991
$ws.Range("J138").Value = 5428.3667
$ws.Range("K138").Value = 32264973
$ws.Range("L138").Value = 16285.1001
$ws.Range("M138").Value = -32259833
$ws.Range("N138").Value = -26565.1001
# Row 141
$ws.Range("H141").Value = 1416.1111
$ws.Range("J141").Value = 2252.5
$ws.Range("L141").Value = 6757.5
$ws.Range("N141").Value = -17117.5

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 2018
$ws.Range("I132").Value = 1551.0312
$ws.Range("J132").Value = 3014.2
$ws.Range("K132").Value = 4653.0936
$ws.Range("L132").Value = 9042.599999999999
$ws.Range("M132").Value = -2123.0936
$ws.Range("N132").Value = -14102.6

$ws = $wb.Worksheets.Item("BSM")
# Row 30
$ws.Range("H30").Value = 1005.5
$ws.Range("J30").Value = 1005.5
$ws.Range("L30").Value = 1005.5
$ws.Range("N30").Value = -1255.5
# Row 55
$ws.Range("H55").Value = 59468
$ws.Range("J55").Value = 59468
$ws.Range("L55").Value = 59468
$ws.Range("N55").Value = -60014
# Row 134
$ws.Range("H134").Value = 2680.0505
$ws.Range("I134").Value = 1774.8679
$ws.Range("J134").Value = 4525.231
$ws.Range("K134").Value = 5324.6037
$ws.Range("L134").Value = 13575.693
$ws.Range("M134").Value = -2789.6037
$ws.Range("N134").Value = -18645.693

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4319.338
$ws.Range("J31").Value = 8017.4414
$ws.Range("L31").Value = 8017.4414
$ws.Range("N31").Value = -8607.4414
# Row 34
$ws.Range("H34").Value = 4319.338
$ws.Range("J34").Value = 8017.4414
$ws.Range("L34").Value = 8017.4414
$ws.Range("N34").Value = -8421.4414
# Row 132
$ws.Range("H132").Value = 3662.4167
$ws.Range("J132").Value = 22532.6
$ws.Range("L132").Value = 67597.79999999999
$ws.Range("N132").Value = -72657.79999999999
# Row 134
$ws.Range("H134").Value = 851998.1
$ws.Range("I134").Value = 1636.36
$ws.Range("K134").Value = 4909.08
$ws.Range("M134").Value = -2374.08

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 172
$ws.Range("I11").Value = 156.15384
$ws.Range("J11").Value = 275
$ws.Range("K11").Value = 468.46152
$ws.Range("L11").Value = 825
$ws.Range("M11").Value = -328.46152
$ws.Range("N11").Value = -1105
# Row 38
$ws.Range("H38").Value = 142.59259
$ws.Range("I38").Value = 161.25
$ws.Range("J38").Value = 127.666664
$ws.Range("K38").Value = 483.75
$ws.Range("L38").Value = 382.999992
$ws.Range("M38").Value = -136.75
$ws.Range("N38").Value = -1076.999992
# Row 107
$ws.Range("H107").Value = 561.4737
$ws.Range("I107").Value = 317.16666
$ws.Range("J107").Value = 674.2308
$ws.Range("K107").Value = 951.4999799999999
$ws.Range("L107").Value = 2022.6924
$ws.Range("M107").Value = 968.5000200000001
$ws.Range("N107").Value = -5862.6924
# Row 118
$ws.Range("H118").Value = 4382.857
$ws.Range("I118").Value = 4693.3335
$ws.Range("J118").Value = 4150
$ws.Range("K118").Value = 14080.0005
$ws.Range("L118").Value = 12450
$ws.Range("M118").Value = -12837.0005
$ws.Range("N118").Value = -14936
# Row 122
$ws.Range("H122").Value = 1374
$ws.Range("I122").Value = 1625.3334
$ws.Range("J122").Value = 1260.9
$ws.Range("K122").Value = 14628.0006
$ws.Range("L122").Value = 11348.1
$ws.Range("M122").Value = -12178.0006
$ws.Range("N122").Value = -16248.1
# Row 132
$ws.Range("H132").Value = 1991.9286
$ws.Range("J132").Value = 2308.3
$ws.Range("L132").Value = 20774.7
$ws.Range("N132").Value = -25834.7

$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 1501.8334
$ws.Range("I13").Value = 400
$ws.Range("J13").Value = 1722.2
$ws.Range("K13").Value = 400
$ws.Range("L13").Value = 1722.2
$ws.Range("M13").Value = -261
$ws.Range("N13").Value = -2000.2
# Row 15
$ws.Range("H15").Value = 19700
$ws.Range("J15").Value = 19700
$ws.Range("L15").Value = 19700
$ws.Range("N15").Value = -20276
# Row 81
$ws.Range("H81").Value = 19700
$ws.Range("J81").Value = 19700
$ws.Range("L81").Value = 19700
$ws.Range("N81").Value = -21696
# Row 84
$ws.Range("H84").Value = 19700
$ws.Range("J84").Value = 19700
$ws.Range("L84").Value = 59100
$ws.Range("N84").Value = -69084
# Row 126
$ws.Range("H126").Value = 3926.8
$ws.Range("I126").Value = 2268.6667
$ws.Range("J126").Value = 5032.222
$ws.Range("K126").Value = 6806.000100000001
$ws.Range("L126").Value = 15096.666
$ws.Range("M126").Value = -4336.000100000001
$ws.Range("N126").Value = -20036.666
# Row 132
$ws.Range("H132").Value = 5170.9375
$ws.Range("I132").Value = 5661.926
$ws.Range("J132").Value = 2519.6
$ws.Range("K132").Value = 16985.778
$ws.Range("L132").Value = 7558.799999999999
$ws.Range("M132").Value = -14455.778
$ws.Range("N132").Value = -12618.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6012.971
$ws.Range("I7").Value = 3761.6365
$ws.Range("J7").Value = 9822.923000000001
$ws.Range("K7").Value = 3761.6365
$ws.Range("L7").Value = 9822.923000000001
$ws.Range("M7").Value = -3649.6365
$ws.Range("N7").Value = -10046.923
# Row 9
$ws.Range("H9").Value = 510
$ws.Range("I9").Value = 510
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 510
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -286
# Row 40
$ws.Range("H40").Value = 3053.8286
$ws.Range("I40").Value = 7087.3335
$ws.Range("J40").Value = 1657.6154
$ws.Range("K40").Value = 7087.3335
$ws.Range("L40").Value = 1657.6154
$ws.Range("M40").Value = -6951.3335
$ws.Range("N40").Value = -1929.6154
# Row 55
$ws.Range("H55").Value = 332.5238
$ws.Range("I55").Value = 220.11111
$ws.Range("J55").Value = 416.83334
$ws.Range("K55").Value = 220.11111
$ws.Range("L55").Value = 416.83334
$ws.Range("M55").Value = -47.11111
$ws.Range("N55").Value = -762.83334
# Row 68
$ws.Range("H68").Value = 2183.15
$ws.Range("I68").Value = 1924.2667
$ws.Range("J68").Value = 2959.8
$ws.Range("K68").Value = 1924.2667
$ws.Range("L68").Value = 2959.8
$ws.Range("M68").Value = -1175.2667
$ws.Range("N68").Value = -4457.8
# Row 71
$ws.Range("H71").Value = 2183.15
$ws.Range("I71").Value = 1924.2667
$ws.Range("J71").Value = 2959.8
$ws.Range("K71").Value = 9621.333499999999
$ws.Range("L71").Value = 14799
$ws.Range("M71").Value = -5877.333499999999
$ws.Range("N71").Value = -22287
# Row 80
$ws.Range("H80").Value = 23929.334
$ws.Range("J80").Value = 23929.334
$ws.Range("L80").Value = 23929.334
$ws.Range("N80").Value = -26175.334
# Row 83
$ws.Range("H83").Value = 23929.334
$ws.Range("J83").Value = 23929.334
$ws.Range("L83").Value = 71788.00199999999
$ws.Range("N83").Value = -83020.00199999999
# Row 126
$ws.Range("H126").Value = 6012.971
$ws.Range("I126").Value = 3761.6365
$ws.Range("J126").Value = 9822.923000000001
$ws.Range("K126").Value = 11284.9095
$ws.Range("L126").Value = 29468.769
$ws.Range("M126").Value = -8814.9095
$ws.Range("N126").Value = -34408.769

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5281.4736
$ws.Range("I62").Value = 5666.6665
$ws.Range("J62").Value = 5103.6924
$ws.Range("K62").Value = 5666.6665
$ws.Range("L62").Value = 5103.6924
$ws.Range("M62").Value = -5042.6665
$ws.Range("N62").Value = -6351.6924
# Row 65
$ws.Range("H65").Value = 5281.4736
$ws.Range("I65").Value = 5666.6665
$ws.Range("J65").Value = 5103.6924
$ws.Range("K65").Value = 28333.3325
$ws.Range("L65").Value = 25518.462
$ws.Range("M65").Value = -25213.3325
$ws.Range("N65").Value = -31758.462
# Row 107
$ws.Range("H107").Value = 959.5
$ws.Range("J107").Value = 292.5
$ws.Range("L107").Value = 877.5
$ws.Range("N107").Value = -4717.5
# Row 132
$ws.Range("H132").Value = 1913.963
$ws.Range("I132").Value = 1745.9756
$ws.Range("J132").Value = 2443.7693
$ws.Range("K132").Value = 5237.9268
$ws.Range("L132").Value = 7331.3079
$ws.Range("M132").Value = -2707.9268
$ws.Range("N132").Value = -12391.3079
# Row 136
$ws.Range("H136").Value = 1307.1818
$ws.Range("I136").Value = 640.4483
$ws.Range("J136").Value = 6141
$ws.Range("K136").Value = 1921.3449
$ws.Range("L136").Value = 18423
$ws.Range("M136").Value = 628.6550999999999
$ws.Range("N136").Value = -23523
